$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# A1, B1, D1 stay the same (Receiver_Name, Receiver_Phone_Number, Receiver_Address)
$ws.Range("C1").Value = "Receiver_City"
$ws.Range("E1").Value = "Receiver_Area"
$ws.Range("F1").Value = "Product_Name"
$ws.Range("G1").Value = "Product_Quantity_Pieces"
$ws.Range("H1").Value = "Product_Cost"
$ws.Range("I1").Value = "Delivery_Note"

# --- Update data row (row 2) ---
# A2, B2 stay the same (Ather, 01982114988)
$ws.Range("D2").Value = "Mirpur Model"
$ws.Range("E2").Value = "Mirpur DOHS, Home: 4/A, Flat No: 5"
$ws.Range("F2").Value = "Chocolate"
$ws.Range("G2").Value = "5"
$ws.Range("H2").Value = "50"
# I2 was previously blank with the default (non-centered) style; give it the same
# centered data-row formatting as the rest of row 2 before putting the value in.
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Value = "Bring Quixx!"
# C2 becomes fully blank (old Receiver_Area value removed / no receiver city data supplied)
$ws.Range("C2").Clear()

# --- New header style for C1: underlined Calibri 11 (minor theme font), centered ---
# Start from the same base formatting used by the other data-row cells (numFmt + centered alignment),
# then layer on the font differences so we end up with a distinct font entry.
$src = $ws.Range("A2")
$src.Copy()
$ws.Range("C1").PasteSpecial(-4122)
$c1Font = $ws.Range("C1").Font
$c1Font.Name = "Calibri"
$c1Font.ThemeFont = 1
$c1Font.Underline = 2

# --- Column width: column G (Product_Quantity_Pieces) widened ---
$ws.Columns.Item(7).ColumnWidth = 29.59

# --- Selection / view state ---
$ws.Range("J1:J2").Select()
$ws.Range("J2").Activate()

Write-Host "edit applied"
